$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) storage for numeric-looking price values so Excel
# does not auto-convert them to numbers (source data is textual).
$textCells = @("D4", "D5", "D6", "D8", "D11", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values.
$ws.Range('D2').Value = '61.530.02'
$ws.Range('E2').Value = '  -1.87%  '
$ws.Range('D3').Value = '3.004.37'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '596.05'
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').Value = '143.84'
$ws.Range('E6').Value = '  -2.62%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '0.522'
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').Value = '3.002.18'
$ws.Range('E9').Value = '  -0.76%  '
$ws.Range('E10').Value = '  -1.77%  '
$ws.Range('D11').Value = '5.89'
$ws.Range('E11').Value = '  +1.68%  '
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  +3.73%  '
$ws.Range('D13').Value = '0.0000229'
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D14').Value = '34.32'
$ws.Range('E14').Value = '  -1.64%  '
$ws.Range('E15').Value = '  +2.26%  '
$ws.Range('D16').Value = '3.497.97'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('D17').Value = '7.03'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D18').Value = '61.532.98'
$ws.Range('E18').Value = '  -1.72%  '
$ws.Range('D19').Value = '3.003.60'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('D20').Value = '453.05'
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').Value = '14.04'
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('D22').Value = '0.687'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').Value = '7.36'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').Value = '81.88'
$ws.Range('E24').Value = '  +2.21%  '
$ws.Range('E25').Value = '  -4.45%  '
$ws.Range('D26').Value = '10.72'
$ws.Range('E26').Value = '  +5.30%  '
$ws.Range('D27').Value = '11.96'
$ws.Range('E27').Value = '  -3.82%  '
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('E29').Value = '  +1.79%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('D31').Value = '7.24'
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('D32').Value = '2.07'
$ws.Range('E32').Value = '  -2.29%  '
$ws.Range('D33').Value = '27.54'
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('E34').Value = '  +1.75%  '
$ws.Range('D35').Value = '0.0₃0841'
$ws.Range('E35').Value = '  +5.42%  '
$ws.Range('D36').Value = '1.02'
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('D37').Value = '5.78'
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').Value = '9.24'
$ws.Range('E38').Value = '  +2.60%  '
$ws.Range('D39').Value = '2.07'
$ws.Range('E39').Value = '  -3.10%  '
$ws.Range('D40').Value = '50.34'
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.124'
$ws.Range('E41').Value = '  +9.47%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '2.90'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('D43').Value = '399.33'
$ws.Range('E43').Value = '  -5.82%  '
$ws.Range('D44').Value = '39.93'
$ws.Range('E44').Value = '  +4.60%  '
$ws.Range('D45').Value = '0.0354'
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('E46').Value = '  -2.55%  '
$ws.Range('D47').Value = '2.715.98'
$ws.Range('E47').Value = '  -2.58%  '
$ws.Range('D48').Value = '132.43'
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('D50').Value = '0.107'
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('E51').Value = '  +1.78%  '
